$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 17:52"

# Update countries & provincias Spain: refreshed case counts and re-sorted rows
$ws.Range("B4").Value = 321615
$ws.Range("C4").Value = 10258
$ws.Range("D4").Value = 16570
$ws.Range("E4").Value = 295913
$ws.Range("G4").Value = 680
$ws.Range("H4").Value = 9132
$ws.Range("B16").Value = 14426
$ws.Range("C16").Value = 514
$ws.Range("E16").Value = 11565
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = 258
$ws.Range("B17").Value = 12051
$ws.Range("C17").Value = 270
$ws.Range("E17").Value = 8849
$ws.Range("E21").Value = 7493
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 48
$ws.Range("E24").Value = 5584
$ws.Range("G24").Value = 8
$ws.Range("H24").Value = 70
$ws.Range("A30").Value = "Polonia"
$ws.Range("B30").Value = 4102
$ws.Range("C30").Value = 475
$ws.Range("D30").Value = 134
$ws.Range("E30").Value = 3874
$ws.Range("F30").Value = 50
$ws.Range("G30").Value = 15
$ws.Range("H30").Value = 94
$ws.Range("A31").Value = "Rumania"
$ws.Range("B31").Value = 3864
$ws.Range("C31").Value = 251
$ws.Range("D31").Value = 374
$ws.Range("E31").Value = 3339
$ws.Range("F31").Value = 141
$ws.Range("H31").Value = 151
$ws.Range("B38").Value = 2804
$ws.Range("C38").Value = 75
$ws.Range("E38").Value = 2268
$ws.Range("F38").Value = 33
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 36
$ws.Range("B48").Value = 1735
$ws.Range("C48").Value = 62
$ws.Range("E48").Value = 1584
$ws.Range("F48").Value = 93
$ws.Range("G48").Value = 5
$ws.Range("H48").Value = 73
$ws.Range("A64").Value = "Irak"
$ws.Range("B64").Value = 961
$ws.Range("C64").Value = 83
$ws.Range("D64").Value = 279
$ws.Range("E64").Value = 621
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 5
$ws.Range("H64").Value = 61
$ws.Range("A65").Value = "Hong Kong"
$ws.Range("B65").Value = 890
$ws.Range("C65").Value = 28
$ws.Range("D65").Value = 206
$ws.Range("E65").Value = 680
$ws.Range("F65").Value = 8
$ws.Range("H65").Value = 4
$ws.Range("D75").Value = 52
$ws.Range("E75").Value = 502
$ws.Range("D85").Value = 37
$ws.Range("E85").Value = 400
$ws.Range("A101").Value = "Mauricio"
$ws.Range("B101").Value = 228
$ws.Range("C101").Value = 11
$ws.Range("D101").Value = 23
$ws.Range("E101").Value = 204
$ws.Range("F101").Value = 0
$ws.Range("H101").Value = 1
$ws.Range("A102").Value = "Estado de Palestina"
$ws.Range("C102").Value = 14
$ws.Range("D102").Value = 2
$ws.Range("E102").Value = 225
$ws.Range("F102").Value = 3
$ws.Range("H102").Value = 0
$ws.Range("A103").Value = "Malta"
$ws.Range("B103").Value = 227
$ws.Range("C103").Value = 31
$ws.Range("D103").Value = 7
$ws.Range("E103").Value = 213
$ws.Range("F103").Value = 1
$ws.Range("H103").Value = 7
$ws.Range("B109").Value = 175
$ws.Range("C109").Value = 9
$ws.Range("D109").Value = 33
$ws.Range("E109").Value = 137
$ws.Range("A122").Value = "Guinea"
$ws.Range("B122").Value = 121
$ws.Range("C122").Value = 10
$ws.Range("D122").Value = 5
$ws.Range("E122").Value = 116
$ws.Range("F122").Value = 0
$ws.Range("A123").Value = "Camboya"
$ws.Range("B123").Value = 114
$ws.Range("D123").Value = 50
$ws.Range("E123").Value = 64
$ws.Range("F123").Value = 1
$ws.Range("A137").Value = "Barbados"
$ws.Range("B137").Value = 56
$ws.Range("C137").Value = 4
$ws.Range("D137").Value = 6
$ws.Range("E137").Value = 49
$ws.Range("F137").Value = 4
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 1
$ws.Range("A138").Value = "Jamaica"
$ws.Range("B138").Value = 55
$ws.Range("C138").Value = 2
$ws.Range("D138").Value = 7
$ws.Range("E138").Value = 45
$ws.Range("H138").Value = 3
$ws.Range("A140").Value = "Macao"
$ws.Range("C140").Value = 4
$ws.Range("D140").Value = 1
$ws.Range("E140").Value = 39
$ws.Range("G140").Value = 2
$ws.Range("A141").Value = "Togo"
$ws.Range("B141").Value = 45
$ws.Range("C141").Value = 23
$ws.Range("D141").Value = 2
$ws.Range("E141").Value = 38
$ws.Range("G141").Value = 3
$ws.Range("H141").Value = 5
$ws.Range("A142").Value = "Etiopia"
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 10
$ws.Range("E142").Value = 34
$ws.Range("H142").Value = 0
$ws.Range("A143").Value = "Mali"
$ws.Range("B143").Value = 44
$ws.Range("C143").Value = 3
$ws.Range("D143").Value = 20
$ws.Range("E143").Value = 21
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 3
$ws.Range("A144").Value = "Congo"
$ws.Range("B144").Value = 43
$ws.Range("C144").Value = 5
$ws.Range("D144").Value = 4
$ws.Range("E144").Value = 38
$ws.Range("F144").Value = 1
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 1
